$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 383.5
$ws.Range("I6").Value = 383.5
$ws.Range("K6").Value = 1150.5
$ws.Range("M6").Value = -1038.5
# Row 17
$ws.Range("H17").Value = 640.2069
$ws.Range("J17").Value = 640.2069
$ws.Range("L17").Value = 1920.6207
$ws.Range("N17").Value = -2256.6207
# Row 20
$ws.Range("H20").Value = 657
$ws.Range("I20").Value = 657
$ws.Range("K20").Value = 657
$ws.Range("M20").Value = -427
# Row 28
$ws.Range("H28").Value = 1044.6
$ws.Range("I28").Value = 466
$ws.Range("J28").Value = 3359
$ws.Range("K28").Value = 466
$ws.Range("L28").Value = 3359
$ws.Range("M28").Value = 19
$ws.Range("N28").Value = -4329
# Row 33
$ws.Range("H33").Value = 5600.8945
$ws.Range("I33").Value = 6810
$ws.Range("J33").Value = 1066.75
$ws.Range("K33").Value = 6810
$ws.Range("L33").Value = 1066.75
$ws.Range("M33").Value = -6581
$ws.Range("N33").Value = -1524.75
# Row 35
$ws.Range("H35").Value = 657
$ws.Range("I35").Value = 657
$ws.Range("K35").Value = 657
$ws.Range("M35").Value = -278
# Row 43
$ws.Range("H43").Value = 4317.5
$ws.Range("I43").Value = 3899.5
$ws.Range("J43").Value = 4363.9443
$ws.Range("K43").Value = 3899.5
$ws.Range("L43").Value = 4363.9443
$ws.Range("M43").Value = -3830.5
$ws.Range("N43").Value = -4501.9443
# Row 53
$ws.Range("H53").Value = 377.1111
$ws.Range("I53").Value = 66.583336
$ws.Range("J53").Value = 998.1667
$ws.Range("K53").Value = 66.583336
$ws.Range("L53").Value = 998.1667
$ws.Range("M53").Value = 570.416664
$ws.Range("N53").Value = -2272.1667
# Row 70
$ws.Range("H70").Value = 2899.1
$ws.Range("I70").Value = 2852
$ws.Range("K70").Value = 8556
$ws.Range("M70").Value = -8286
# Row 73
$ws.Range("H73").Value = 2899.1
$ws.Range("I73").Value = 2852
$ws.Range("K73").Value = 8556
$ws.Range("M73").Value = -7620
# Row 80
$ws.Range("H80").Value = 1140.174
$ws.Range("I80").Value = 1505.8
$ws.Range("J80").Value = 1038.6111
$ws.Range("K80").Value = 4517.4
$ws.Range("L80").Value = 3115.8333
$ws.Range("M80").Value = -3519.4
$ws.Range("N80").Value = -5111.8333
# Row 83
$ws.Range("H83").Value = 1140.174
$ws.Range("I83").Value = 1505.8
$ws.Range("J83").Value = 1038.6111
$ws.Range("K83").Value = 13552.2
$ws.Range("L83").Value = 9347.4999
$ws.Range("M83").Value = -8560.199999999999
$ws.Range("N83").Value = -19331.4999
# Row 101
$ws.Range("H101").Value = 547.5
$ws.Range("I101").Value = 495
$ws.Range("K101").Value = 1485
$ws.Range("M101").Value = 137
# Row 127
$ws.Range("H127").Value = 1856.5714
$ws.Range("I127").Value = 1830
$ws.Range("K127").Value = 5490
$ws.Range("M127").Value = -530
# Row 129
$ws.Range("H129").Value = 1136.1923
$ws.Range("I129").Value = 954.1
$ws.Range("K129").Value = 2862.3
$ws.Range("M129").Value = 2137.7
# Row 141
$ws.Range("H141").Value = 1309.5
$ws.Range("I141").Value = 1345.3334
$ws.Range("J141").Value = 987
$ws.Range("K141").Value = 4036.0002
$ws.Range("L141").Value = 2961
$ws.Range("M141").Value = 1143.9998
$ws.Range("N141").Value = -13321

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2666.3333
$ws.Range("I2").Value = 1999.5
$ws.Range("K2").Value = 1999.5
$ws.Range("M2").Value = -1886.5
# Row 32
$ws.Range("H32").Value = 28042.139
$ws.Range("I32").Value = 28946.77
$ws.Range("K32").Value = 28946.77
$ws.Range("M32").Value = -28659.77
# Row 45
$ws.Range("H45").Value = 3827.9167
$ws.Range("I45").Value = 2871.2144
$ws.Range("J45").Value = 5167.3
$ws.Range("K45").Value = 2871.2144
$ws.Range("L45").Value = 5167.3
$ws.Range("M45").Value = -2494.2144
$ws.Range("N45").Value = -5921.3
# Row 61
$ws.Range("H61").Value = 6676.5454
$ws.Range("I61").Value = 3710
$ws.Range("J61").Value = 9148.667
$ws.Range("K61").Value = 3710
$ws.Range("L61").Value = 9148.667
$ws.Range("M61").Value = -3498
$ws.Range("N61").Value = -9572.667
# Row 74
$ws.Range("H74").Value = 917.25
$ws.Range("I74").Value = 897.6
$ws.Range("K74").Value = 897.6
$ws.Range("M74").Value = -23.60000000000002
# Row 77
$ws.Range("H77").Value = 917.25
$ws.Range("I77").Value = 897.6
$ws.Range("K77").Value = 4488
$ws.Range("M77").Value = -120
# Row 88
$ws.Range("H88").Value = 3438.1
$ws.Range("J88").Value = 3567.5
$ws.Range("L88").Value = 3567.5
$ws.Range("N88").Value = -4379.5
# Row 91
$ws.Range("H91").Value = 3438.1
$ws.Range("J91").Value = 3567.5
$ws.Range("L91").Value = 3567.5
$ws.Range("N91").Value = -6375.5
# Row 110
$ws.Range("H110").Value = 2581.75
$ws.Range("I110").Value = 2080
$ws.Range("J110").Value = 3083.5
$ws.Range("K110").Value = 2080
$ws.Range("L110").Value = 3083.5
$ws.Range("M110").Value = -35
$ws.Range("N110").Value = -7173.5
# Row 112
$ws.Range("H112").Value = 18999.5
$ws.Range("J112").Value = 18999.5
$ws.Range("L112").Value = 18999.5
$ws.Range("N112").Value = -21953.5
# Row 116
$ws.Range("H116").Value = 2666.3333
$ws.Range("I116").Value = 1999.5
$ws.Range("K116").Value = 1999.5
$ws.Range("M116").Value = 294.5
# Row 122
$ws.Range("H122").Value = 2535.5386
$ws.Range("I122").Value = 2346.8333
$ws.Range("K122").Value = 7040.499899999999
$ws.Range("M122").Value = -4590.499899999999
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 136
$ws.Range("H136").Value = 6676.5454
$ws.Range("I136").Value = 3710
$ws.Range("J136").Value = 9148.667
$ws.Range("K136").Value = 11130
$ws.Range("L136").Value = 27446.001
$ws.Range("M136").Value = -8580
$ws.Range("N136").Value = -32546.001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2666.3333
$ws.Range("I3").Value = 1999.5
$ws.Range("K3").Value = 1999.5
$ws.Range("M3").Value = -1885.5
# Row 20
$ws.Range("H20").Value = 5559250.5
$ws.Range("I20").Value = 10002582
$ws.Range("J20").Value = 5086.75
$ws.Range("K20").Value = 10002582
$ws.Range("L20").Value = 5086.75
$ws.Range("M20").Value = -10002335
$ws.Range("N20").Value = -5580.75
# Row 35
$ws.Range("H35").Value = 44000.5
$ws.Range("J35").Value = 44000.5
$ws.Range("L35").Value = 44000.5
$ws.Range("N35").Value = -44620.5
# Row 86
$ws.Range("H86").Value = 3123.0322
$ws.Range("I86").Value = 1777.3334
$ws.Range("J86").Value = 4384.625
$ws.Range("K86").Value = 1777.3334
$ws.Range("L86").Value = 4384.625
$ws.Range("M86").Value = -654.3334
$ws.Range("N86").Value = -6630.625
# Row 89
$ws.Range("H89").Value = 3123.0322
$ws.Range("I89").Value = 1777.3334
$ws.Range("J89").Value = 4384.625
$ws.Range("K89").Value = 8886.667
$ws.Range("L89").Value = 21923.125
$ws.Range("M89").Value = -3270.666999999999
$ws.Range("N89").Value = -33155.125
# Row 134
$ws.Range("H134").Value = 2554.2246
$ws.Range("I134").Value = 1967.3334
$ws.Range("J134").Value = 6075.5713
$ws.Range("K134").Value = 5902.0002
$ws.Range("L134").Value = 18226.7139
$ws.Range("M134").Value = -3367.0002
$ws.Range("N134").Value = -23296.7139

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 4666.3
$ws.Range("I33").Value = 1804.1428
$ws.Range("J33").Value = 11344.667
$ws.Range("K33").Value = 1804.1428
$ws.Range("L33").Value = 11344.667
$ws.Range("M33").Value = -1425.1428
$ws.Range("N33").Value = -12102.667
# Row 58
$ws.Range("H58").Value = 74144.21
$ws.Range("I58").Value = 93445.91
$ws.Range("J58").Value = 3371.3333
$ws.Range("K58").Value = 93445.91
$ws.Range("L58").Value = 3371.3333
$ws.Range("M58").Value = -93242.91
$ws.Range("N58").Value = -3777.3333
# Row 62
$ws.Range("H62").Value = 4486.7144
$ws.Range("I62").Value = 4533
$ws.Range("J62").Value = 4452
$ws.Range("K62").Value = 4533
$ws.Range("L62").Value = 4452
$ws.Range("M62").Value = -3909
$ws.Range("N62").Value = -5700
# Row 65
$ws.Range("H65").Value = 4486.7144
$ws.Range("I65").Value = 4533
$ws.Range("J65").Value = 4452
$ws.Range("K65").Value = 22665
$ws.Range("L65").Value = 22260
$ws.Range("M65").Value = -19545
$ws.Range("N65").Value = -28500
# Row 105
$ws.Range("H105").Value = 1090
$ws.Range("I105").Value = 1090
$ws.Range("K105").Value = 1090
$ws.Range("M105").Value = 657
# Row 107
$ws.Range("H107").Value = 2346.6177
$ws.Range("I107").Value = 740.3333
$ws.Range("J107").Value = 3614.7368
$ws.Range("K107").Value = 740.3333
$ws.Range("L107").Value = 3614.7368
$ws.Range("M107").Value = 1179.6667
$ws.Range("N107").Value = -7454.736800000001
# Row 136
$ws.Range("H136").Value = 74144.21
$ws.Range("I136").Value = 93445.91
$ws.Range("J136").Value = 3371.3333
$ws.Range("K136").Value = 280337.73
$ws.Range("L136").Value = 10113.9999
$ws.Range("M136").Value = -277787.73
$ws.Range("N136").Value = -15213.9999

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 591928.4
$ws.Range("I4").Value = 418421.88
$ws.Range("K4").Value = 1255265.64
$ws.Range("M4").Value = -1255153.64
# Row 9
$ws.Range("H9").Value = 695
$ws.Range("I9").Value = 695
$ws.Range("K9").Value = 2085
$ws.Range("M9").Value = -1861
# Row 23
$ws.Range("H23").Value = 436
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 537.75
$ws.Range("K23").Value = 87
$ws.Range("L23").Value = 1613.25
$ws.Range("M23").Value = 148
$ws.Range("N23").Value = -2083.25
# Row 37
$ws.Range("H37").Value = 63858.3
$ws.Range("J37").Value = 63858.3
$ws.Range("L37").Value = 191574.9
$ws.Range("N37").Value = -191798.9
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 48214.215
# Row 80
$ws.Range("H80").Value = 3431.72
$ws.Range("I80").Value = 2111.4546
$ws.Range("K80").Value = 2111.4546
$ws.Range("M80").Value = -1113.4546
# Row 83
$ws.Range("H83").Value = 3431.72
$ws.Range("I83").Value = 2111.4546
$ws.Range("K83").Value = 10557.273
$ws.Range("M83").Value = -5565.273000000001
# Row 97
$ws.Range("H97").Value = 2610.6667
$ws.Range("I97").Value = 2499.4285
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2499.4285
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -2003.4285
$ws.Range("N97").Value = -3992
# Row 107
$ws.Range("H107").Value = 48456
$ws.Range("J107").Value = 1196.4546
$ws.Range("L107").Value = 1196.4546
$ws.Range("N107").Value = -5036.4546
# Row 122
$ws.Range("H122").Value = 4913.7144
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850
# Row 132
$ws.Range("H132").Value = 55199.105
$ws.Range("I132").Value = 65099.125
$ws.Range("J132").Value = 2399
$ws.Range("K132").Value = 195297.375
$ws.Range("L132").Value = 7197
$ws.Range("M132").Value = -192767.375
$ws.Range("N132").Value = -12257

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1912.0741
$ws.Range("I16").Value = 1735
$ws.Range("J16").Value = 2930.25
$ws.Range("K16").Value = 1735
$ws.Range("L16").Value = 2930.25
$ws.Range("M16").Value = -1565
$ws.Range("N16").Value = -3270.25
# Row 38
$ws.Range("H38").Value = 11126.4
$ws.Range("I38").Value = 3500
$ws.Range("K38").Value = 3500
$ws.Range("M38").Value = -3090
# Row 46
$ws.Range("H46").Value = 4473.7705
$ws.Range("J46").Value = 2844.7551
$ws.Range("L46").Value = 2844.7551
$ws.Range("N46").Value = -3220.7551
# Row 68
$ws.Range("H68").Value = 4436.3
$ws.Range("I68").Value = 2802.6924
$ws.Range("J68").Value = 7470.143
$ws.Range("K68").Value = 2802.6924
$ws.Range("L68").Value = 7470.143
$ws.Range("M68").Value = -2053.6924
$ws.Range("N68").Value = -8968.143
# Row 71
$ws.Range("H71").Value = 4436.3
$ws.Range("I71").Value = 2802.6924
$ws.Range("J71").Value = 7470.143
$ws.Range("K71").Value = 14013.462
$ws.Range("L71").Value = 37350.715
$ws.Range("M71").Value = -10269.462
$ws.Range("N71").Value = -44838.715
# Row 82
$ws.Range("H82").Value = 1695.963
$ws.Range("I82").Value = 1321.8948
$ws.Range("K82").Value = 1321.8948
$ws.Range("M82").Value = -960.8948
# Row 85
$ws.Range("H85").Value = 1695.963
$ws.Range("I85").Value = 1321.8948
$ws.Range("K85").Value = 1321.8948
$ws.Range("M85").Value = -73.89480000000003
# Row 132
$ws.Range("H132").Value = 41764.324
$ws.Range("I132").Value = 48730.348
$ws.Range("J132").Value = 5541
$ws.Range("K132").Value = 146191.044
$ws.Range("L132").Value = 16623
$ws.Range("M132").Value = -143661.044
$ws.Range("N132").Value = -21683

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 867.7586
$ws.Range("I107").Value = 472.85715
$ws.Range("K107").Value = 1418.57145
$ws.Range("M107").Value = 501.4285500000001
# Row 122
$ws.Range("H122").Value = 929.8182
$ws.Range("I122").Value = 828.17645
$ws.Range("K122").Value = 2484.52935
$ws.Range("M122").Value = -34.52935000000025
